$d = $word.ActiveDocument

# Update the date heading (first paragraph) via Find/Replace to avoid touching paragraph marks
$d.Content.Find.Execute("2023-11-05 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-06 Monday", 2) | Out-Null

# Update the table cells (20 rows x 5 cols), in row-major order matching the diff
$t = $d.Tables.Item(1)
$values = @(
    "4+44=48",
    "68-2=66",
    "21+55=76",
    "17+29=46",
    "82-0=82",
    "29+30=59",
    "19+36=55",
    "36-8=28",
    "54-24=30",
    "42+17=59",
    "42+1=43",
    "28+52=80",
    "13+22=35",
    "41+3=44",
    "11+7=18",
    "91-11=80",
    "67-62=5",
    "39-30=9",
    "85-59=26",
    "26+29=55",
    "54-27=27",
    "47+18=65",
    "24+28=52",
    "28+51=79",
    "81-56=25",
    "72-36=36",
    "84-8=76",
    "17+28=45",
    "16+56=72",
    "67-0=67",
    "41+2=43",
    "83-43=40",
    "96-57=39",
    "12-3=9",
    "97-50=47",
    "1+46=47",
    "73-39=34",
    "41+39=80",
    "97-47=50",
    "92-69=23",
    "2+9=11",
    "7+19=26",
    "64-1=63",
    "9+47=56",
    "1+72=73",
    "8+20=28",
    "73-66=7",
    "79-74=5",
    "17+7=24",
    "4+68=72",
    "67-67=0",
    "4+57=61",
    "79-9=70",
    "64+2=66",
    "14+83=97",
    "62-50=12",
    "82-27=55",
    "52+46=98",
    "38-37=1",
    "88-52=36",
    "31-5=26",
    "67+31=98",
    "67-49=18",
    "14+2=16",
    "75-64=11",
    "64+13=77",
    "22+46=68",
    "49+1=50",
    "97-46=51",
    "27-8=19",
    "27+40=67",
    "21+21=42",
    "61-7=54",
    "86-75=11",
    "79-35=44",
    "16+40=56",
    "96-69=27",
    "55-25=30",
    "82-64=18",
    "77-12=65",
    "89-58=31",
    "20+25=45",
    "85-85=0",
    "52-1=51",
    "47-1=46",
    "69-13=56",
    "9+19=28",
    "47-37=10",
    "94-4=90",
    "56-50=6",
    "74+12=86",
    "17+19=36",
    "66+3=69",
    "29+4=33",
    "89-28=61",
    "97-33=64",
    "61+19=80",
    "47+38=85",
    "99-1=98",
    "96-11=85"
)

$numCols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [int][Math]::Floor($i / $numCols) + 1
    $col = ($i % $numCols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}

Write-Host "Done"
